# Applies the NIH-NCPI ncpi-Study-Participant StructureDefinition regeneration
# (FHIR version downgraded from 4.3.0/R4B to 4.0.1/R4, IG regenerated on a new date)
# to the workbook that is already open as $excel.ActiveWorkbook.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsEl   = $wb.Worksheets.Item("Elements")

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------

# Date the IG was generated
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version: 4.3.0 (R4B) -> 4.0.1 (R4)
$wsMeta.Range("B15").Value = "4.0.1"

# ---------------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------------

# ResearchSubject (row 2) - Constraint(s): dom-3 invariant text simplified
# (drops the "id.exists() and" clause and fixes a copy/paste duplicate clause)
$aj2 = @'
dom-2:If the resource is contained in another resource, it SHALL NOT contain nested Resources {contained.contained.empty()}
dom-3:If the resource is contained in another resource, it SHALL be referred to from elsewhere in the resource or SHALL refer to the containing resource {contained.where((('#'+id in (%resource.descendants().reference | %resource.descendants().as(canonical) | %resource.descendants().as(uri) | %resource.descendants().as(url))) or descendants().where(reference = '#').exists() or descendants().where(as(canonical) = '#').exists() or descendants().where(as(canonical) = '#').exists()).not()).trace('unmatched', id).empty()}dom-4:If a resource is contained in another resource, it SHALL NOT have a meta.versionId or a meta.lastUpdated {contained.meta.versionId.empty() and contained.meta.lastUpdated.empty()}dom-5:If a resource is contained in another resource, it SHALL NOT have a security label {contained.meta.security.empty()}dom-6:A resource should have narrative for robust management {text.`div`.exists()}
'@
$wsEl.Range("AJ2").Value = $aj2

# ResearchSubject (row 2) - Mapping: FiveWs Pattern Mapping was mis-populated
# with an HL7 v2 style value ("clinical.diagnostics"); clear it.
$wsEl.Range("AN2").Value = ""

# ResearchSubject.language (row 6) - Binding Description
$wsEl.Range("Y6").Value = "A human language."

# ResearchSubject.contained (row 8) - Constraint(s): the dom-r4b invariant only
# applies to R4B; now that this IG targets FHIR 4.0.1 (R4), drop it entirely.
$wsEl.Range("AJ8").Value = ""

# DomainResource.modifierExtension (row 11) - Comments: fix hyperlink from the
# R4B spec to the R4 spec
$o11 = @'
Modifier extensions allow for extensions that *cannot* be safely ignored to be clearly distinguished from the vast majority of extensions which can be safely ignored.  This promotes interoperability by eliminating the need for implementers to prohibit the presence of extensions. For further information, see the [definition of modifier extensions](http://hl7.org/fhir/R4/extensibility.html#modifierExtension).
'@
$wsEl.Range("O11").Value = $o11

# ResearchSubject.extension:access-policy (row 10) - Constraint(s): align with
# the plain ele-1/ext-1 text used elsewhere (drop the R4B-only "unless an empty
# Parameters resource" wording) - same text as AJ9/AJ11
$aj10 = @'
ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}
'@
$wsEl.Range("AJ10").Value = $aj10

# ResearchSubject.status (row 13) - Binding Value Set
$wsEl.Range("Z13").Value = "http://hl7.org/fhir/ValueSet/research-subject-status|4.0.1"
